$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.627.46'
$ws.Range("E2").Value = '  +4.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.135.19'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.90'
$ws.Range("E5").Value = '  +1.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '610.40'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  +1.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.383'
$ws.Range("E8").Value = '  -2.54%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.131.88'
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.779'
$ws.Range("E11").Value = '  -7.52%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '97.223.36'
$ws.Range("E13").Value = '  +4.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000238'
$ws.Range("E14").Value = '  -2.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.44'
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '33.68'
$ws.Range("E16").Value = '  -4.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.714.71'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.129.60'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '516.23'
$ws.Range("E19").Value = '  +16.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.42'
$ws.Range("E20").Value = '  -9.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.37'
$ws.Range("E21").Value = '  -3.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.61'
$ws.Range("E22").Value = '  -7.35%  '
$ws.Range("E23").Value = '  -5.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.69'
$ws.Range("E24").Value = '  -4.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.16'
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.43'
$ws.Range("E26").Value = '  -4.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.49'
$ws.Range("E27").Value = '  -11.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.292.79'
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -2.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '8.88'
$ws.Range("E34").Value = '  -4.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.52'
$ws.Range("E35").Value = '  +1.94%  '
$ws.Range("E36").Value = '  -6.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.22'
$ws.Range("E37").Value = '  -10.23%  '
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.28'
$ws.Range("E39").Value = '  +1.14%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.432'
$ws.Range("E40").Value = '  -3.98%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '463.01'
$ws.Range("E41").Value = '  -3.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.21'
$ws.Range("E42").Value = '  -7.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.52'
$ws.Range("E43").Value = '  -10.66%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.08'
$ws.Range("E45").Value = '  -7.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '162.31'
$ws.Range("E46").Value = '  +2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.693'
$ws.Range("E47").Value = '  -1.72%  '
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.45'
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.09'
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("E51").Value = '  +0.00%  '
